$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the receipt number (P3) -> new blank entry
$ws.Range("P3").Value = ""

# Update the date (N5)
$ws.Range("N5").Value = "12 de junio de 2023"

# Update the recipient name / DNI line (L7)
$ws.Range("L7").Value = "ALVAREZ MARIA ALISA , DNI 33051226"

# Update the concept (C8) and amount (E8) for the first installment row
$ws.Range("C8").Value = "Cuota1"
$ws.Range("E8").Value = 1

# Update the address-ish field (L9) - must stay text, not become a number,
# so stage it through a helper cell (as a formula returning text) and paste
# only the resulting value back, which keeps the cell's existing style/format.
$ws.Range("ZZ1").Formula = "=""20"""
$ws.Range("ZZ1").Copy()
$ws.Range("L9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# Update the amount in words (J11)
$ws.Range("J11").Value = "UN PESO "

# Update the concept description (I13) - this cell's style uses a
# quote-prefix text format, so stage it through a helper cell too (the same
# formula + paste-values trick), to avoid Excel silently dropping the
# quotePrefix formatting on a plain assignment.
$ws.Range("ZZ1").Formula = "=""En concepto de pago en efectivo por Cuota1"""
$ws.Range("ZZ1").Copy()
$ws.Range("I13").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# Update the sub-total (P15)
$ws.Range("P15").Value = 1
